$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 42500
$ws.Range("J3").Value = 42500
$ws.Range("L3").Value = 42500
$ws.Range("N3").Value = -42728
$ws.Range("H32").Value = 1048.15
$ws.Range("I32").Value = 818.2
$ws.Range("J32").Value = 1278.1
$ws.Range("K32").Value = 818.2
$ws.Range("L32").Value = 1278.1
$ws.Range("M32").Value = -492.2
$ws.Range("N32").Value = -1930.1
$ws.Range("H70").Value = 2128.4285
$ws.Range("I70").Value = 887.375
$ws.Range("J70").Value = 3783.1667
$ws.Range("K70").Value = 2662.125
$ws.Range("L70").Value = 11349.5001
$ws.Range("M70").Value = -2392.125
$ws.Range("N70").Value = -11889.5001
$ws.Range("H73").Value = 2128.4285
$ws.Range("I73").Value = 887.375
$ws.Range("J73").Value = 3783.1667
$ws.Range("K73").Value = 2662.125
$ws.Range("L73").Value = 11349.5001
$ws.Range("M73").Value = -1726.125
$ws.Range("N73").Value = -13221.5001
$ws.Range("H100").Value = 5000
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459
$ws.Range("H102").Value = 42500
$ws.Range("J102").Value = 42500
$ws.Range("L102").Value = 42500
$ws.Range("N102").Value = -48990
$ws.Range("H132").Value = 4332.467
$ws.Range("I132").Value = 4534.7856
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 13604.3568
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -11074.3568
$ws.Range("N132").Value = -9560
$ws.Range("H135").Value = 643.4
$ws.Range("I135").Value = 277.33334
$ws.Range("K135").Value = 2496.00006
$ws.Range("M135").Value = 38.9999399999997
$ws.Range("H137").Value = 2963.4119
$ws.Range("I137").Value = 1357
$ws.Range("J137").Value = 4087.9
$ws.Range("K137").Value = 4071
$ws.Range("L137").Value = 12263.7
$ws.Range("M137").Value = -1521
$ws.Range("N137").Value = -17363.7

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 736.3333
$ws.Range("I30").Value = 604.5
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 604.5
$ws.Range("L30").Value = 1000
$ws.Range("M30").Value = -454.5
$ws.Range("N30").Value = -1300
$ws.Range("H32").Value = 11338
$ws.Range("I32").Value = 12042
$ws.Range("K32").Value = 12042
$ws.Range("M32").Value = -11755
$ws.Range("H45").Value = 1982.4
$ws.Range("J45").Value = 2100
$ws.Range("L45").Value = 2100
$ws.Range("N45").Value = -2854
$ws.Range("H61").Value = 3666.3333
$ws.Range("I61").Value = 3499.5
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3499.5
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3287.5
$ws.Range("N61").Value = -4424
$ws.Range("H132").Value = 1635.3334
$ws.Range("I132").Value = 1635.3334
$ws.Range("K132").Value = 4906.0002
$ws.Range("M132").Value = -2376.0002
$ws.Range("H136").Value = 3666.3333
$ws.Range("I136").Value = 3499.5
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 10498.5
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -7948.5
$ws.Range("N136").Value = -17100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 6000
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H20").Value = 5222.3335
$ws.Range("I20").Value = 4712.923
$ws.Range("K20").Value = 4712.923
$ws.Range("M20").Value = -4465.923
$ws.Range("H49").Value = 17500
$ws.Range("J49").Value = 17500
$ws.Range("L49").Value = 17500
$ws.Range("N49").Value = -17978
$ws.Range("H86").Value = 3500
$ws.Range("I86").Value = 3500
$ws.Range("K86").Value = 3500
$ws.Range("M86").Value = -2377
$ws.Range("H89").Value = 3500
$ws.Range("I89").Value = 3500
$ws.Range("K89").Value = 17500
$ws.Range("M89").Value = -11884
$ws.Range("H99").Value = 1496.3334
$ws.Range("I99").Value = 1496.3334
$ws.Range("K99").Value = 1496.3334
$ws.Range("M99").Value = 1.666600000000017
$ws.Range("H107").Value = 1094.5
$ws.Range("I107").Value = 792.8333
$ws.Range("J107").Value = 1999.5
$ws.Range("K107").Value = 792.8333
$ws.Range("L107").Value = 1999.5
$ws.Range("M107").Value = 1127.1667
$ws.Range("N107").Value = -5839.5
$ws.Range("H134").Value = 3619.0833
$ws.Range("I134").Value = 3734.5454
$ws.Range("J134").Value = 2349
$ws.Range("K134").Value = 11203.6362
$ws.Range("L134").Value = 7047
$ws.Range("M134").Value = -8668.636200000001
$ws.Range("N134").Value = -12117
$ws.Range("H137").Value = 46999.2
$ws.Range("J137").Value = 49999
$ws.Range("L137").Value = 49999
$ws.Range("N137").Value = -60199

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 1817.375
$ws.Range("I25").Value = 910
$ws.Range("J25").Value = 2724.75
$ws.Range("K25").Value = 910
$ws.Range("L25").Value = 2724.75
$ws.Range("M25").Value = -736
$ws.Range("N25").Value = -3072.75
$ws.Range("H31").Value = 2526.3416
$ws.Range("I31").Value = 2047.5714
$ws.Range("K31").Value = 2047.5714
$ws.Range("M31").Value = -1752.5714
$ws.Range("H34").Value = 2526.3416
$ws.Range("I34").Value = 2047.5714
$ws.Range("K34").Value = 2047.5714
$ws.Range("M34").Value = -1845.5714
$ws.Range("H95").Value = 16872.273
$ws.Range("J95").Value = 16872.273
$ws.Range("L95").Value = 16872.273
$ws.Range("N95").Value = -22364.273
$ws.Range("H106").Value = 21835
$ws.Range("J106").Value = 21835
$ws.Range("L106").Value = 21835
$ws.Range("N106").Value = -24359
$ws.Range("H107").Value = 958.1667
$ws.Range("I107").Value = 416.33334
$ws.Range("K107").Value = 416.33334
$ws.Range("M107").Value = 1503.66666
$ws.Range("H132").Value = 1047.125
$ws.Range("I132").Value = 1047.125
$ws.Range("K132").Value = 3141.375
$ws.Range("M132").Value = -611.375

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1233.6666
$ws.Range("I8").Value = 1233.6666
$ws.Range("K8").Value = 3700.9998
$ws.Range("M8").Value = -3561.9998
$ws.Range("H37").Value = 98333.336
$ws.Range("J37").Value = 98333.336
$ws.Range("L37").Value = 295000.008
$ws.Range("N37").Value = -295224.008
$ws.Range("H68").Value = 1529.375
$ws.Range("I68").Value = 1559
$ws.Range("K68").Value = 4677
$ws.Range("M68").Value = -3866
$ws.Range("H71").Value = 1529.375
$ws.Range("I71").Value = 1559
$ws.Range("K71").Value = 14031
$ws.Range("M71").Value = -9975
$ws.Range("H131").Value = 1418
$ws.Range("I131").Value = 999.1667
$ws.Range("J131").Value = 1920.6
$ws.Range("K131").Value = 2997.5001
$ws.Range("L131").Value = 5761.799999999999
$ws.Range("M131").Value = 2042.4999
$ws.Range("N131").Value = -15841.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2856.6667
$ws.Range("I80").Value = 1575
$ws.Range("J80").Value = 3497.5
$ws.Range("K80").Value = 1575
$ws.Range("L80").Value = 3497.5
$ws.Range("M80").Value = -577
$ws.Range("N80").Value = -5493.5
$ws.Range("H83").Value = 2856.6667
$ws.Range("I83").Value = 1575
$ws.Range("J83").Value = 3497.5
$ws.Range("K83").Value = 7875
$ws.Range("L83").Value = 17487.5
$ws.Range("M83").Value = -2883
$ws.Range("N83").Value = -27471.5
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5189.5
$ws.Range("I132").Value = 5741.6
$ws.Range("J132").Value = 4977.154
$ws.Range("K132").Value = 17224.8
$ws.Range("L132").Value = 14931.462
$ws.Range("M132").Value = -14694.8
$ws.Range("N132").Value = -19991.462

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11998.5
$ws.Range("I136").Value = 11998.5
$ws.Range("K136").Value = 35995.5
$ws.Range("M136").Value = -33445.5
